$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Blank000"
$ws.Range("B3").Value = "Blank001"
